# Daily update: advance the tracking date by one day for every active row.
# Column E = "剩余" (days remaining), column F = "开始时间" (cycle start date, yyyymmdd int).
# Row 36 has a malformed 9-digit date and is skipped by the source process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99
$skipRows = @(36)

for ($row = 2; $row -le $lastRow; $row++) {
    if ($skipRows -contains $row) {
        continue
    }

    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $e = $eCell.Value2

    if ($e -eq 1) {
        # Cycle completed -> reset remaining days to 10 and roll start date forward 10 days.
        $eCell.Value2 = 10
        $fCell.Value2 = $fCell.Value2 + 10
    } else {
        $eCell.Value2 = $e - 1
    }
}
